$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Update firstName and lastName for row 2
$ws.Range("B2").Value = "Noah"
$ws.Range("C2").Value = "Seligson"

# Remove "Technology, " from the interests list in G2
$ws.Range("G2").Value = "Education, Environment, Sports & Recreation, Coding & Software Development, Music & Performance, Health & Wellness, Animal Welfare"
